$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 69
$ws1.Range("F4").Value = 1760
$ws1.Range("G4").Value = 60
$ws1.Range("F7").Value = 1143
$ws1.Range("F8").Value = 1577
$ws1.Range("F9").Value = 171
$ws1.Range("F10").Value = 171
$ws1.Range("F11").Value = 19
$ws1.Range("F12").Value = 1500
$ws1.Range("F13").Value = 3125
$ws1.Range("F14").Value = 669
$ws1.Range("F15").Value = 1812
$ws1.Range("F16").Value = 0
$ws1.Range("F17").Value = 879
$ws1.Range("F18").Value = 296
$ws1.Range("F20").Value = 1505
$ws1.Range("F21").Value = 302
$ws1.Range("F24").Value = 1276
$ws1.Range("F27").Value = 173
$ws1.Range("F28").Value = 5839
$ws1.Range("F29").Value = 5338
$ws1.Range("F30").Value = 767
$ws1.Range("F31").Value = 590
$ws1.Range("F32").Value = 1696
$ws1.Range("F33").Value = 88
$ws1.Range("F34").Value = 218

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 32

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 69
$ws4.Range("F7").Value = 1760
$ws4.Range("G7").Value = 60
$ws4.Range("F10").Value = 1143
$ws4.Range("F11").Value = 1577
$ws4.Range("F12").Value = 171
$ws4.Range("F13").Value = 171
$ws4.Range("F15").Value = 19
$ws4.Range("F16").Value = 1500
$ws4.Range("F17").Value = 3125
$ws4.Range("F18").Value = 669
$ws4.Range("F19").Value = 1812
$ws4.Range("F20").Value = 1821
$ws4.Range("F21").Value = 879
$ws4.Range("F22").Value = 296
$ws4.Range("F24").Value = 1505
$ws4.Range("F25").Value = 302
$ws4.Range("F29").Value = 32
$ws4.Range("F30").Value = 1276
$ws4.Range("F33").Value = 173
$ws4.Range("F34").Value = 5839
$ws4.Range("F35").Value = 5338
$ws4.Range("F36").Value = 767
$ws4.Range("F37").Value = 590
$ws4.Range("F38").Value = 1696
$ws4.Range("F41").Value = 88
$ws4.Range("F42").Value = 218
